# Repair tests after last changes in cost methods.
#
# This script reproduces, via the Excel COM object model, the set of edits
# described by the target diff:
#  - mobility!main_min_density value 700000 -> 800000, and its description
#    text updated to refer to "net density".
#  - mobility's old min_weekly_freq row (row 26) is replaced in place by a
#    new mobility_cost_rpc shadow-price-ratio parameter (value 0.82, with a
#    dark-gray font).
#  - infrastructure gains a new infrast_cost_rpc shadow-price-ratio
#    parameter row (value 0.82, dark-gray font).
#  - time gains the min_weekly_freq row (moved here from mobility) plus a
#    new time_cost_rpc shadow-price-ratio parameter row (value 1, dark-gray
#    font).
#  - the "time" sheet becomes the active/selected sheet/tab.

$wb = $excel.ActiveWorkbook

$mobility = $wb.Worksheets.Item("mobility")
$infrastructure = $wb.Worksheets.Item("infrastructure")
$time = $wb.Worksheets.Item("time")

# ---------------------------------------------------------------------
# 1. mobility sheet
# ---------------------------------------------------------------------

# main_min_density: bump value and refresh the description wording.
$mobility.Range("B24").Value = 800000
$mobility.Range("C24").Value = "Minimum net density to consider a link as being a main track (ton-km/ton = ton)."

# Row 26 used to hold min_weekly_freq; that parameter moves to the "time"
# sheet (re-created below) and this row becomes mobility_cost_rpc instead.
$mobility.Range("A26").Value = "mobility_cost_rpc"
$mobility.Range("B26").Style = "Normal"
$mobility.Range("B26").Value = 0.82
$mobility.Range("B26").Font.Color = 2236962
$mobility.Range("C26").Value = "Shadow price to market price ratio in mobility cost (coeff)."

# ---------------------------------------------------------------------
# 2. infrastructure sheet - append infrast_cost_rpc
# ---------------------------------------------------------------------

$infrastructure.Range("A15").Value = "infrast_cost_rpc"
$infrastructure.Range("B15").Value = 0.82
$infrastructure.Range("B15").Font.Color = 2236962
$infrastructure.Range("C15").Value = "Shadow price to market price ratio in infrastructure cost (coeff)."

# ---------------------------------------------------------------------
# 3. time sheet - re-add min_weekly_freq (moved from mobility) and append
#    the new time_cost_rpc row
# ---------------------------------------------------------------------

$time.Range("A6").Value = "min_weekly_freq"
$time.Range("B6").Value = 2
$time.Range("B6").HorizontalAlignment = -4108
$time.Range("C6").Value = "Minimum trains per week that have to be to service an od pair (number)."

$time.Range("A7").Value = "time_cost_rpc"
$time.Range("B7").Value = 1
$time.Range("B7").Font.Color = 2236962
$time.Range("C7").Value = "Shadow price to market price ratio in time cost (coeff)."

# ---------------------------------------------------------------------
# 4. Selections / active sheet. mobility and infrastructure keep their own
#    per-sheet selection state; "time" is activated last so it ends up the
#    active tab.
# ---------------------------------------------------------------------

$mobility.Activate()
$mobility.Range("B27").Select()

$infrastructure.Activate()
$infrastructure.Range("A15:C15").Select()

$time.Activate()
$time.Range("C7").Select()

Write-Output "ok"
